# "Arquivos de entrada atualizados"
# The only substantive content change is the parameter label in row 2,
# column F: "Turb" (abbreviation) is renamed to "Turbidez" (full word).
# Selecting F3 afterwards mirrors the author's recorded cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("F2").Value = "Turbidez"

$ws.Range("F3").Select()
